# Rename the "Device" sheet to "Apparatus" and update its contents to
# reflect the new "apparatus" terminology (was "device"), per the commit:
# Change "Device" to "Apparatus" in excel form, simulink, function name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Device")

# Remove the old "Notes:" row (row 2) — the new sheet no longer has it,
# which shifts every row below it up by one.
[void]$ws.Rows.Item(2).Delete()

# Update the sheet's description text (row 1) to mention "apparatuses"
# instead of "devices".
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."

# The header row (now row 3) column labels change from
# "Device type" / "Device parameters" to "Type" / "Parameters".
$ws.Range("B3").Value = "Type"
$ws.Range("C3").Value = "Parameters"

# Update the selected cell to match the new layout.
[void]$ws.Range("C4").Select()

# Finally, rename the worksheet tab itself.
$ws.Name = "Apparatus"
